# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 15e0a490-... entry (row 3) on both the zh-cn and
# de-de report sheets to reflect the newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 04:42:30"
$wsZhCn.Range("H3").Value = "2016-03-22 04:42:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 04:42:34"
$wsDeDe.Range("H3").Value = "2016-03-22 04:43:00"
